$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing user row (row 2): reorder name fields, new mail ---
$ws.Range("A2").Value = "iva"
$ws.Range("B2").Value = "ivanka"
$ws.Range("C2").Value = "horoshko"
$ws.Range("E2").Value = "aaaa@gmail.com"

# --- Add new user row (row 3) ---
$ws.Range("A3").Value = "myroslava"
$ws.Range("B3").Value = "myroslav"
$ws.Range("C3").Value = "shram"
$ws.Range("D3").Value = "qwerty"
$ws.Range("E3").Value = "gmail@gmail.com"
$ws.Range("F3").Value = "Administrator"

# --- Rewrite header row (row 1) with descriptive column names ---
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Mail"
$ws.Range("F1").Value = "Role"
$ws.Range("B1").Value = "FirstName"
$ws.Range("A1").Value = "Login Name"

# --- Widen columns A:F (~12.71 characters wide) ---
for ($col = 1; $col -le 6; $col++) {
    $ws.Columns.Item($col).ColumnWidth = 11.8
}

# --- Move active selection to B1 ---
[void]$ws.Range("B1").Select()
